$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("14/04/2022", "13:45", "13:52", 7, 3, 4, 2, 2, "EURUSD",     -0.5600000000000001, "LOSS"),
    @("16/04/2022", "16:00", "16:07", 7, 3, 4, 2, 1, "EURUSD-OTC", -0.5500000000000003, "LOSS"),
    @("16/04/2022", "16:26", "16:30", 3, 0, 3, 3, 0, "EURUSD-OTC", -9.98,               "LOSS"),
    @("16/04/2022", "16:30", "16:38", 7, 3, 4, 3, 2, "EURUSD-OTC", -5.82,               "LOSS"),
    @("17/04/2022", "13:45", "13:56", 10, 6, 4, 1, 2, "EURUSD-OTC", 5.789999999999999,  "LOSS")
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]

    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
}
